$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 830.3333
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H62").Value = 76958184
$ws.Range("I62").Value = 142858430
$ws.Range("K62").Value = 142858430
$ws.Range("M62").Value = -142857806
$ws.Range("H65").Value = 76958184
$ws.Range("I65").Value = 142858430
$ws.Range("K65").Value = 714292150
$ws.Range("M65").Value = -714289030
$ws.Range("H123").Value = 99380
$ws.Range("J123").Value = 99380
$ws.Range("L123").Value = 99380
$ws.Range("N123").Value = -109180
$ws.Range("H127").Value = 3566
$ws.Range("I127").Value = 1421.3334
$ws.Range("K127").Value = 4264.0002
$ws.Range("M127").Value = 695.9997999999996
$ws.Range("H137").Value = 3844.4546
$ws.Range("I137").Value = 5900
$ws.Range("K137").Value = 17700
$ws.Range("M137").Value = -15150
$ws.Range("H138").Value = 5201.41
$ws.Range("I138").Value = 1465
$ws.Range("J138").Value = 9560.556
$ws.Range("K138").Value = 4395
$ws.Range("L138").Value = 28681.668
$ws.Range("M138").Value = 745
$ws.Range("N138").Value = -38961.66800000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1959.4375
$ws.Range("I2").Value = 1139.9412
$ws.Range("J2").Value = 2888.2
$ws.Range("K2").Value = 1139.9412
$ws.Range("L2").Value = 2888.2
$ws.Range("M2").Value = -1026.9412
$ws.Range("N2").Value = -3114.2
$ws.Range("H16").Value = 817.5
$ws.Range("I16").Value = 1022.5
$ws.Range("K16").Value = 1022.5
$ws.Range("M16").Value = -735.5
$ws.Range("H34").Value = 273799.8
$ws.Range("I34").Value = 84500
$ws.Range("J34").Value = 399999.66
$ws.Range("K34").Value = 84500
$ws.Range("L34").Value = 399999.66
$ws.Range("M34").Value = -84229
$ws.Range("N34").Value = -400541.66
$ws.Range("H61").Value = 5080.62
$ws.Range("I61").Value = 2416.8
$ws.Range("K61").Value = 2416.8
$ws.Range("M61").Value = -2204.8
$ws.Range("H116").Value = 1959.4375
$ws.Range("I116").Value = 1139.9412
$ws.Range("J116").Value = 2888.2
$ws.Range("K116").Value = 1139.9412
$ws.Range("L116").Value = 2888.2
$ws.Range("M116").Value = 1154.0588
$ws.Range("N116").Value = -7476.2
$ws.Range("H132").Value = 4760.1406
$ws.Range("I132").Value = 3545.9773
$ws.Range("K132").Value = 10637.9319
$ws.Range("M132").Value = -8107.9319
$ws.Range("H136").Value = 5080.62
$ws.Range("I136").Value = 2416.8
$ws.Range("K136").Value = 7250.400000000001
$ws.Range("M136").Value = -4700.400000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1959.4375
$ws.Range("I3").Value = 1139.9412
$ws.Range("J3").Value = 2888.2
$ws.Range("K3").Value = 1139.9412
$ws.Range("L3").Value = 2888.2
$ws.Range("M3").Value = -1025.9412
$ws.Range("N3").Value = -3116.2
$ws.Range("H20").Value = 15153450
$ws.Range("I20").Value = 20835494
$ws.Range("J20").Value = 1332.3334
$ws.Range("K20").Value = 20835494
$ws.Range("L20").Value = 1332.3334
$ws.Range("M20").Value = -20835247
$ws.Range("N20").Value = -1826.3334
$ws.Range("H87").Value = 74998
$ws.Range("J87").Value = 74998
$ws.Range("L87").Value = 74998
$ws.Range("N87").Value = -77494
$ws.Range("H90").Value = 74998
$ws.Range("J90").Value = 74998
$ws.Range("L90").Value = 224994
$ws.Range("N90").Value = -237474
$ws.Range("H107").Value = 35159844
$ws.Range("I107").Value = 45003228
$ws.Range("J107").Value = 4898.5713
$ws.Range("K107").Value = 45003228
$ws.Range("L107").Value = 4898.5713
$ws.Range("M107").Value = -45001308
$ws.Range("N107").Value = -8738.5713
$ws.Range("H134").Value = 5483.9375
$ws.Range("I134").Value = 2120.1538
$ws.Range("K134").Value = 6360.4614
$ws.Range("M134").Value = -3825.4614
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13165253
$ws.Range("I58").Value = 33336672
$ws.Range("J58").Value = 9980.869000000001
$ws.Range("K58").Value = 33336672
$ws.Range("L58").Value = 9980.869000000001
$ws.Range("M58").Value = -33336469
$ws.Range("N58").Value = -10386.869
$ws.Range("H59").Value = 75025.75
$ws.Range("I59").Value = 104
$ws.Range("J59").Value = 99999.664
$ws.Range("K59").Value = 104
$ws.Range("L59").Value = 99999.664
$ws.Range("M59").Value = 1041
$ws.Range("N59").Value = -102289.664
$ws.Range("H64").Value = 69997.5
$ws.Range("J64").Value = 69997.5
$ws.Range("L64").Value = 69997.5
$ws.Range("N64").Value = -70493.5
$ws.Range("H67").Value = 69997.5
$ws.Range("J67").Value = 69997.5
$ws.Range("L67").Value = 69997.5
$ws.Range("N67").Value = -71713.5
$ws.Range("H95").Value = 47876.375
$ws.Range("J95").Value = 47876.375
$ws.Range("L95").Value = 47876.375
$ws.Range("N95").Value = -53368.375
$ws.Range("H132").Value = 5032.7856
$ws.Range("I132").Value = 2231.4092
$ws.Range("J132").Value = 8114.3
$ws.Range("K132").Value = 6694.2276
$ws.Range("L132").Value = 24342.9
$ws.Range("M132").Value = -4164.2276
$ws.Range("N132").Value = -29402.9
$ws.Range("H136").Value = 13165253
$ws.Range("I136").Value = 33336672
$ws.Range("J136").Value = 9980.869000000001
$ws.Range("K136").Value = 100010016
$ws.Range("L136").Value = 29942.607
$ws.Range("M136").Value = -100007466
$ws.Range("N136").Value = -35042.607
$ws.Range("H141").Value = 74048
$ws.Range("J141").Value = 74048
$ws.Range("L141").Value = 74048
$ws.Range("N141").Value = -84408
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 180.91667
$ws.Range("J23").Value = 154.83333
$ws.Range("L23").Value = 464.49999
$ws.Range("N23").Value = -934.49999
$ws.Range("H134").Value = 96944.27
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H141").Value = 7287.9287
$ws.Range("I141").Value = 2504
$ws.Range("J141").Value = 13666.5
$ws.Range("K141").Value = 7512
$ws.Range("L141").Value = 40999.5
$ws.Range("M141").Value = -2332
$ws.Range("N141").Value = -51359.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H70").Value = 8802.536
$ws.Range("I70").Value = 7721.8335
$ws.Range("K70").Value = 7721.8335
$ws.Range("M70").Value = -7451.8335
$ws.Range("H73").Value = 8802.536
$ws.Range("I73").Value = 7721.8335
$ws.Range("K73").Value = 7721.8335
$ws.Range("M73").Value = -6785.8335
$ws.Range("H97").Value = 4647.3
$ws.Range("I97").Value = 3781.8572
$ws.Range("K97").Value = 3781.8572
$ws.Range("M97").Value = -3285.8572
$ws.Range("H132").Value = 6860.4287
$ws.Range("I132").Value = 2592.2856
$ws.Range("J132").Value = 15396.714
$ws.Range("K132").Value = 7776.8568
$ws.Range("L132").Value = 46190.142
$ws.Range("M132").Value = -5246.8568
$ws.Range("N132").Value = -51250.142
$ws.Range("H133").Value = 71365
$ws.Range("J133").Value = 71365
$ws.Range("L133").Value = 71365
$ws.Range("N133").Value = -81485
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 35000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -4710
$ws.Range("N33").ClearContents()
$ws.Range("H46").Value = 1838
$ws.Range("I46").Value = 1056
$ws.Range("K46").Value = 1056
$ws.Range("M46").Value = -868
$ws.Range("H122").Value = 4640.3335
$ws.Range("I122").Value = 3915.8262
$ws.Range("J122").Value = 5922.154
$ws.Range("K122").Value = 11747.4786
$ws.Range("L122").Value = 17766.462
$ws.Range("M122").Value = -9297.4786
$ws.Range("N122").Value = -22666.462
$ws.Range("H136").Value = 9380.107
$ws.Range("I136").Value = 2428.111
$ws.Range("J136").Value = 12673.158
$ws.Range("K136").Value = 7284.333
$ws.Range("L136").Value = 38019.474
$ws.Range("M136").Value = -4734.333
$ws.Range("N136").Value = -43119.474
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 21002
$ws.Range("I15").Value = 21002
$ws.Range("K15").Value = 21002
$ws.Range("M15").Value = -20714
$ws.Range("H136").Value = 27058492
$ws.Range("I136").Value = 71429550
$ws.Range("K136").Value = 214288650
$ws.Range("M136").Value = -214286100
